# Generate Report for Handoff
# Replace the old localization GUID/hash-based file names and refresh the
# handoff/handback timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "7c1907f1-4c53-4a18-8761-caac5578edaa"
$newGuid = "1cb91d82-c797-4adc-bdaf-aeab64b649ef"

$oldHash = "5e166d4e1991a47ad7e6db05065d95b373ff6174"
$newHash = "f90b542ae74e50539d3d5afef1203f0252100bff"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38d279af640c65ab3b7ecd97cd33fdad8b370c9b/e2e/$oldGuid.md"

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 15:07:04"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-22 15:06:56"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-22 15:07:04"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
